# Update the "Last Updated" timestamp on the Metadata sheet
$wb = $excel.ActiveWorkbook
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("A2").Value = "29 Oct 2025, 07:07 PM"

# Insert a new worksheet named "1 Month Performance" right before
# the existing "distance from Dma50" sheet.
$targetSheet = $wb.Worksheets.Item("distance from Dma50")
$newSheet = $wb.Worksheets.Add($targetSheet)
$newSheet.Name = "1 Month Performance"

# Header row
$newSheet.Cells.Item(1, 1).Value = "Icon"
$newSheet.Cells.Item(1, 2).Value = "Stock"
$newSheet.Cells.Item(1, 3).Value = "% Change"

$headerRange = $newSheet.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows (Icon, Stock, % Change)
$data = @(
    @("INDOTHAI", 106.9036),
    @("VISASTEEL", 79.9766),
    @("GMBREW", 77.5378),
    @("MCLEODRUSS", 69.8655),
    @("PROZONER", 68.3711),
    @("IFBAGRO", 66.1892),
    @("BGRENERGY", 64.929),
    @("ESSARSHPNG", 64.91160000000001),
    @("PANACHE", 62.6487),
    @("MAHASTEEL", 55.9703),
    @("INOXGREEN", 51.0181),
    @("STALLION", 46.4325),
    @("ORIENTTECH", 45.3321),
    @("TVSSRICHAK", 40.7778),
    @("MTARTECH", 40.7213),
    @("SEJALLTD", 37.4301),
    @("V2RETAIL", 37.2004),
    @("RAMAPHO", 36.9731),
    @("SANDUMA", 36.9057),
    @("TARACHAND", 36.4813),
    @("NETWEB", 36.1199),
    @("SAMMAANCAP", 35.5128),
    @("ONMOBILE", 35.4702),
    @("SHAREINDIA", 35.3207),
    @("SOUTHBANK", 35.2819),
    @("TVSELECT", 35.1983),
    @("RAMCOSYS", 34.6928),
    @("MAANALU", 34.4803),
    @("MEGASOFT", 33.4399),
    @("BHARATSE", 31.8611),
    @("EMKAY", 30.3743),
    @("ATHERENERG", 29.116),
    @("TATVA", 28.7451),
    @("TERASOFT", 28.3093),
    @("CARTRADE", 27.5713),
    @("ARFIN", 27.3801),
    @("MINDTECK", 26.9415),
    @("BHARATWIRE", 26.5276),
    @("HATSUN", 26.492),
    @("INDORAMA", 26.4516),
    @("IFBIND", 26.161),
    @("ADANIPOWER", 25.8247),
    @("AVALON", 25.7352),
    @("MRPL", 25.6265),
    @("HINDCOPPER", 25.3164),
    @("PRECWIRE", 24.679),
    @("SCI", 24.132),
    @("KICL", 24.1119),
    @("SKYGOLD", 24.1079),
    @("DCBBANK", 23.8922),
    @("AUBANK", 23.6964),
    @("ETHOSLTD", 23.1527),
    @("PVP", 22.7524),
    @("INDIANB", 22.6689),
    @("PRIVISCL", 22.3984),
    @("CPEDU", 22.3786),
    @("LORDSCHLO", 22.1791),
    @("GUJTHEM", 22.0704),
    @("SURYODAY", 21.8039),
    @("TDPOWERSYS", 21.7743),
    @("ORBTEXP", 21.6115),
    @("CEATLTD", 20.0239),
    @("ATL", 19.9362),
    @("GRMOVER", 19.7859),
    @("FEDERALBNK", 19.6872),
    @("SUBROS", 19.6508),
    @("USHAMART", 19.6172),
    @("BANKINDIA", 19.3067),
    @("RBLBANK", 19.2556),
    @("MOLDTECH", 19.1891),
    @("THOMASCOTT", 19.1649),
    @("IIFL", 18.9853),
    @("KARURVYSYA", 18.8614),
    @("LUMAXIND", 18.8057),
    @("REPRO", 18.689)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = "📈"
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $r = $r + 1
}

Write-Output "1 Month Performance sheet populated with $($data.Count) rows"
